# Append the 2025-04-07 price row to every "Solar_Prices" sheet.
# Each sheet has dates in column A and a single price value in column B,
# with the new row simply repeating the prior day's (2025-04-06) price.

$wb = $excel.ActiveWorkbook

$newDate = "2025-04-07"

$sheetValues = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "42"
    "N-type Wafer"               = "1.28"
    "Cell Topcon 183mm"          = "0.303"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,487"
    "Silver Busbar front-side"   = "8,215"
    "Silver finger front-side"   = "8,265"
    "USD_CNY"                    = "7.3068"
}

foreach ($sheetName in $sheetValues.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find the first empty row right after the existing data (row 37).
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $newRow = $lastRow + 1

    $dateCell = $ws.Cells.Item($newRow, 1)
    $valueCell = $ws.Cells.Item($newRow, 2)

    # Force plain text so "2025-04-07" and values like "5,487" aren't
    # reinterpreted as a real date / thousands-formatted number.
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $newDate

    $valueCell.NumberFormat = "@"
    $valueCell.Value = $sheetValues[$sheetName]
}
